# ICAS: End of year run 2024
# Updates modelcoefficient results (offset "C" column and "gewijzigd" timestamp
# "E" column) for each Filterweerstand sheet to match the refreshed model run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("IK106")
$ws.Range("E2").Value = 45096.49076744213
$ws.Range("E3").Value = 45096.49076744213
$ws.Range("E4").Value = 45096.49076744213
$ws.Range("E5").Value = 45096.49076744213

$ws = $wb.Worksheets.Item("Q100")
$ws.Range("C2").Value = -0.0200142947273236
$ws.Range("E2").Value = 45657.57032079861
$ws.Range("C3").Value = -0.0692546947273236
$ws.Range("E3").Value = 45657.5703208912
$ws.Range("C4").Value = -0.04835544815926658
$ws.Range("E4").Value = 45657.57032079861
$ws.Range("C5").Value = -0.01629934320602976
$ws.Range("E5").Value = 45657.57032079861

$ws = $wb.Worksheets.Item("Q200")
$ws.Range("C2").Value = -0.04136553704688054
$ws.Range("E2").Value = 45657.57041466435
$ws.Range("C3").Value = [double]"-4.940656458412465e-324"
$ws.Range("E3").Value = 45657.57041466435
$ws.Range("C4").Value = [double]"-6.564079436993284e-147"
$ws.Range("E4").Value = 45657.57041466435

$ws = $wb.Worksheets.Item("Q300")
$ws.Range("C2").Value = -0.0866734723037329
$ws.Range("E2").Value = 45657.57057532408
$ws.Range("C3").Value = -0.0918874723037329
$ws.Range("E3").Value = 45657.57057546296
$ws.Range("C4").Value = -0.09841847230373291
$ws.Range("E4").Value = 45657.57057546296
$ws.Range("C5").Value = -0.1054744723037329
$ws.Range("E5").Value = 45657.57057546296

$ws = $wb.Worksheets.Item("Q400")
$ws.Range("C2").Value = -0.06923117033577117
$ws.Range("E2").Value = 45657.57073831018
$ws.Range("C3").Value = -0.03393798164471695
$ws.Range("E3").Value = 45657.57073831018
$ws.Range("C4").Value = -0.04036048164471694
$ws.Range("E4").Value = 45657.57073841435
$ws.Range("C5").Value = -0.04484048164471695
$ws.Range("E5").Value = 45657.57073841435

$ws = $wb.Worksheets.Item("Q500")
$ws.Range("C2").Value = -0.007381263530535781
$ws.Range("E2").Value = 45657.57097861111
$ws.Range("C3").Value = -0.01373876353053578
$ws.Range("E3").Value = 45657.57097871528
$ws.Range("C4").Value = -0.01916376353053578
$ws.Range("E4").Value = 45657.57097871528
$ws.Range("C5").Value = -0.008998939823900305
$ws.Range("E5").Value = 45657.57097861111
$ws.Range("C6").Value = [double]"-1e-08"
$ws.Range("E6").Value = 45657.57097861111

$ws = $wb.Worksheets.Item("Q600")
$ws.Range("C2").Value = [double]"-4.940656458412465e-324"
$ws.Range("E2").Value = 45657.57106283565
$ws.Range("C3").Value = [double]"-4.940656458412465e-324"
$ws.Range("E3").Value = 45657.57106283565
$ws.Range("C4").Value = -0.1805246804431886
$ws.Range("E4").Value = 45657.57106283565
$ws.Range("C5").Value = [double]"-8.08634922390439e-174"
$ws.Range("E5").Value = 45657.57106283565

$ws = $wb.Worksheets.Item("P100")
$ws.Range("C2").Value = [double]"-8.478591186539953e-97"
$ws.Range("E2").Value = 45657.57116181713
$ws.Range("C3").Value = -0.0600543450924556
$ws.Range("E3").Value = 45657.57116181713
$ws.Range("C4").Value = -0.0003186787480661763
$ws.Range("E4").Value = 45657.57116181713
$ws.Range("C5").Value = -0.04902615100419187
$ws.Range("E5").Value = 45657.57116181713

$ws = $wb.Worksheets.Item("P200")
$ws.Range("C2").Value = [double]"-1.922848456476216e-30"
$ws.Range("E2").Value = 45657.57132618056
$ws.Range("C3").Value = [double]"-8.874685183736383e-29"
$ws.Range("E3").Value = 45657.57132618056
$ws.Range("E4").Value = 45657.57132628472
$ws.Range("C5").Value = [double]"-3.3526588471893e-30"
$ws.Range("E5").Value = 45657.57132618056
$ws.Range("C6").Value = -0.00385
$ws.Range("E6").Value = 45657.57132628472

$ws = $wb.Worksheets.Item("P300")
$ws.Range("C2").Value = [double]"-1.148933078226618e-72"
$ws.Range("E2").Value = 45657.57153850694
$ws.Range("C3").Value = [double]"-7.373037232174746e-77"
$ws.Range("E3").Value = 45657.57153850694
$ws.Range("C4").Value = [double]"-4.890623964652136e-79"
$ws.Range("E4").Value = 45657.57153850694
$ws.Range("C5").Value = -0.2300149801531476
$ws.Range("E5").Value = 45657.57153850694
$ws.Range("C6").Value = -0.2336519743029264
$ws.Range("E6").Value = 45657.57153850694
$ws.Range("C7").Value = -0.1376703730996233
$ws.Range("E7").Value = 45657.57153850694
$ws.Range("C8").Value = -0.08124475998823612
$ws.Range("E8").Value = 45657.57153850694

$ws = $wb.Worksheets.Item("P400")
$ws.Range("C2").Value = -0.04086481189390601
$ws.Range("E2").Value = 45657.57169365741
$ws.Range("C3").Value = -0.03873412725521125
$ws.Range("E3").Value = 45657.57169365741
$ws.Range("C4").Value = -0.04412412725521125
$ws.Range("E4").Value = 45657.57169380787
$ws.Range("C5").Value = -0.0164932564953982
$ws.Range("E5").Value = 45657.57169365741
$ws.Range("C6").Value = [double]"-4.674920645069725e-174"
$ws.Range("E6").Value = 45657.57169365741

$ws = $wb.Worksheets.Item("P500")
$ws.Range("C2").Value = -0.02256892312566483
$ws.Range("E2").Value = 45657.57188424769
$ws.Range("C3").Value = -0.02525560165766922
$ws.Range("E3").Value = 45657.57188424769
$ws.Range("C4").Value = -0.01772351951977796
$ws.Range("E4").Value = 45657.57188424769
$ws.Range("C5").Value = -0.02235086951977796
$ws.Range("E5").Value = 45657.57188435186

$ws = $wb.Worksheets.Item("P600")
$ws.Range("C2").Value = -0.3542975794951402
$ws.Range("E2").Value = 45657.57196103978
$ws.Range("C3").Value = -0.1426303880304105
$ws.Range("E3").Value = 45657.57196103978
$ws.Range("C4").Value = [double]"-1.02827751205387e-62"
$ws.Range("E4").Value = 45657.57196103978
$ws.Range("C5").Value = [double]"-3.802005000714391e-57"
$ws.Range("E5").Value = 45657.57196103978
